$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.11186715098494915,
    -0.0059999999512179159,
    -0.0039999999599551472,
    -0.0079999999256745724,
    -0.0029999999631797891,
    -0.0019999999664666035,
    -0.0099999999029911635,
    -0.0099999999026398889,
    -0.0019999999674555902,
    -0.0019999999689641612,
    -0.0029999999611192152,
    -0.0034999999574387708,
    -0.0034999999611855515,
    -0.0079999999270254918,
    -0.0009999999849474861,
    -0.0019999999800077717,
    -0.0019999999862205797,
    0.031570077527657503,
    -0.0039999999671214148,
    0.014663107535112019,
    -0.0039999999675419673,
    -0.0039999999673012709,
    -0.004999999949172107,
    -0.019999999825015102,
    -0.019999999822530867,
    -0.0024999999624402136,
    -0.0024999999618571245,
    -0.0019999999626971743,
    -0.0069999999208292252,
    -0.059999999498919987,
    -0.0069999999272027935,
    0.03263561508637558,
    0.02218949294190864
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
